$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the stimuli paths to point to the new "stimuli" folder
$ws.Range("A2").Value = "stimuli/duck_video.mp4"
$ws.Range("A3").Value = "stimuli/P07s.mp4"

# Move the active selection to C14 as recorded in the saved view
$ws.Activate()
$ws.Range("C14").Select()
